$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New vm_pu values for rows 2-25, columns B,C,D,E,F,I,J,K,L,M
$rows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25)
$Bvals = @(1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02,1.02)
$Cvals = @(1.050419571049633,1.051887087344584,1.052835602104131,1.053234108588681,1.053301005151543,1.052840927938493,1.050915749588486,1.047514894375465,1.04524158976055,1.04425569726348,1.043889254436922,1.043967868580332,1.044225411893696,1.044384061003105,1.045306987176702,1.045885497640717,1.046222785263936,1.046337766796285,1.04582344427979,1.044149578450846,1.043095769134394,1.043654547074526,1.045851483979362,1.048395137739023)
$Dvals = @(1.047334886536162,1.048394878751425,1.049079622230744,1.049367218352626,1.049415491271675,1.049083466161004,1.047693354483374,1.045234895650314,1.043589716594611,1.042875814162632,1.042610404919305,1.04266734676674,1.042853880152634,1.042968778377017,1.043637063360245,1.044055848429578,1.044299971442547,1.044383186239452,1.04401093201573,1.042798957217119,1.042035583974374,1.042440392630184,1.044031228265008,1.045871543919771)
$Evals = @(1.05661466016842,1.057921602330407,1.05876656877249,1.059121626713343,1.059181232906591,1.058771313720415,1.057056497718347,1.054029107566314,1.052006745230569,1.051130001098873,1.050804176093073,1.050874074140893,1.051103071664931,1.051244142798901,1.052064909282095,1.052579469804299,1.052879503536771,1.052981790376487,1.052524272789261,1.051035642153273,1.050098733890086,1.050595498395286,1.052549214255971,1.054812461920956)
$Fvals = @(1.065336006971663,1.06673091686561,1.067632869009195,1.068011899566101,1.068075531805007,1.067637934217777,1.065807560436635,1.062577023048455,1.060419530944579,1.059484340965872,1.059136816112562,1.059211368426105,1.059455617623475,1.059606087068953,1.060481575215136,1.061030478635776,1.06135055041825,1.061459670816172,1.060971596291837,1.059383696614035,1.058384428969623,1.058914246285722,1.060998202980229,1.063412841091094)
$Ivals = @(1.038470868084027,1.038778445179996,1.038975630911724,1.03905808930625,1.039071908759386,1.038976734446124,1.038575196424498,1.037853503502916,1.037362787392534,1.037148008293797,1.037067883107864,1.037085085947829,1.037141392202978,1.037176038355975,1.037376993055359,1.037502430751057,1.037575375046074,1.037600209666509,1.037488995385638,1.037124820998022,1.036893843529409,1.03701647984357,1.037495066933158,1.038041761922937)
$Jvals = @(1.055452823498969,1.056567485010936,1.057287264884252,1.057589510242754,1.057640238165003,1.057291304871773,1.055829837794783,1.05324300936773,1.051510423445233,1.050758226867596,1.050478525326045,1.050538536021186,1.050735112848226,1.050856190105555,1.051560302169916,1.052001440345493,1.052258558549573,1.052346197110386,1.05195413011016,1.050677234259528,1.049872647756082,1.050299342127415,1.051975508157149,1.053913162087667)
$Kvals = @(1.050097812781029,1.05096898477225,1.051530953772042,1.051766792862529,1.05180636717928,1.05153410668299,1.050392591188771,1.048367639406063,1.047008415385327,1.046417613701292,1.046197821646438,1.04624498333183,1.046399452622426,1.046494580806666,1.047047577215144,1.047393852280127,1.047595611982072,1.047664370198069,1.047356722685737,1.046353974740938,1.045721526496468,1.046056988431739,1.047373500621019,1.048892754000793)
$Lvals = @(1.059351886840596,1.06047124802195,1.061194290505423,1.061497960383681,1.06154893058831,1.06119834931968,1.059730443444195,1.057133972381478,1.055396111794411,1.054641906787763,1.05436150078925,1.054421660709295,1.054618733688918,1.054740122154742,1.055446129631349,1.055888530635109,1.056146411999361,1.056234315280426,1.055841082160398,1.054560707846963,1.053754174850071,1.054181878236642,1.055862522589706,1.057806415224203)
$Mvals = @(1.068049522327009,1.069258363041623,1.07003938436981,1.070367447874068,1.070422515016368,1.070043769056651,1.068458302633665,1.065655272752547,1.063780091702193,1.062966513118142,1.06266406589046,1.062728953156727,1.062941517825973,1.063072452914332,1.06383405177526,1.064311346874478,1.064589589799085,1.064684437284763,1.064260153724488,1.06287892971788,1.062009062186423,1.062470333533347,1.064283286185358,1.066381046848879)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $rows[$i]
    $ws.Range("B" + $r).Value = $Bvals[$i]
    $ws.Range("C" + $r).Value = $Cvals[$i]
    $ws.Range("D" + $r).Value = $Dvals[$i]
    $ws.Range("E" + $r).Value = $Evals[$i]
    $ws.Range("F" + $r).Value = $Fvals[$i]
    $ws.Range("I" + $r).Value = $Ivals[$i]
    $ws.Range("J" + $r).Value = $Jvals[$i]
    $ws.Range("K" + $r).Value = $Kvals[$i]
    $ws.Range("L" + $r).Value = $Lvals[$i]
    $ws.Range("M" + $r).Value = $Mvals[$i]
}

$wb.Save()